$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the price cells that look like plain numbers stay stored as text
# (matching the original inline-string "Price" column formatting) instead of
# being auto-coerced into floating point numbers by Excel.
$textCells = @('D5','D6','D8','D13','D14','D15','D16','D19','D20','D22','D24','D26','D28','D30','D31','D32','D34','D41','D43','D48','D49')
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range('D2').Value = '37.290.74'
$ws.Range('E2').Value = '  +0.12%  '
$ws.Range('D3').Value = '2.062.72'
$ws.Range('E3').Value = '  -0.26%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '233.79'
$ws.Range('E5').Value = '  -0.92%  '
$ws.Range('D6').Value = '0.624'
$ws.Range('E6').Value = '  +1.08%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').Value = '56.57'
$ws.Range('E8').Value = '  -1.17%  '
$ws.Range('E9').Value = '  +0.44%  '
$ws.Range('E10').Value = '  +0.09%  '
$ws.Range('E11').Value = '  +0.57%  '
$ws.Range('D12').Value = '2.365.66'
$ws.Range('E12').Value = '  -0.52%  '
$ws.Range('D13').Value = '14.60'
$ws.Range('E13').Value = '  +0.63%  '
$ws.Range('D14').Value = '20.57'
$ws.Range('E14').Value = '  -3.12%  '
$ws.Range('D15').Value = '0.775'
$ws.Range('E15').Value = '  +0.08%  '
$ws.Range('D16').Value = '5.13'
$ws.Range('E16').Value = '  -2.21%  '
$ws.Range('D17').Value = '2.064.19'
$ws.Range('E17').Value = '  +0.36%  '
$ws.Range('D18').Value = '37.293.02'
$ws.Range('E18').Value = '  -0.45%  '
$ws.Range('D19').Value = '6.33'
$ws.Range('E19').Value = '  +6.84%  '
$ws.Range('D20').Value = '69.22'
$ws.Range('E20').Value = '  +1.32%  '
$ws.Range('D21').Value = '0.0₃0807'
$ws.Range('E21').Value = '  -0.34%  '
$ws.Range('D22').Value = '225.57'
$ws.Range('E22').Value = '  +0.79%  '
$ws.Range('D24').Value = '2.43'
$ws.Range('E24').Value = '  +0.82%  '
$ws.Range('E25').Value = '  -2.15%  '
$ws.Range('D26').Value = '166.19'
$ws.Range('E26').Value = '  +1.80%  '
$ws.Range('E27').Value = '  +4.76%  '
$ws.Range('D28').Value = '8.74'
$ws.Range('E28').Value = '  -1.35%  '
$ws.Range('E29').Value = '  -3.05%  '
$ws.Range('D30').Value = '18.98'
$ws.Range('E30').Value = '  -1.33%  '
$ws.Range('D31').Value = '0.118'
$ws.Range('E31').Value = '  -0.60%  '
$ws.Range('D32').Value = '4.44'
$ws.Range('E32').Value = '  -0.61%  '
$ws.Range('E33').Value = '  -1.10%  '
$ws.Range('D34').Value = '4.53'
$ws.Range('E34').Value = '  +3.54%  '
$ws.Range('E35').Value = '  -1.45%  '
$ws.Range('E36').Value = '  -0.07%  '
$ws.Range('E37').Value = '  -1.28%  '
$ws.Range('E38').Value = '  -3.19%  '
$ws.Range('E39').Value = '  -4.36%  '
$ws.Range('E40').Value = '  -0.32%  '
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').Value = '96.04'
$ws.Range('E41').Value = '  +1.61%  '
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').Value = '1.460.34'
$ws.Range('E42').Value = '  -0.70%  '
$ws.Range('D43').Value = '0.0933'
$ws.Range('E43').Value = '  -2.29%  '
$ws.Range('E44').Value = '  +1.68%  '
$ws.Range('E45').Value = '  +2.07%  '
$ws.Range('E46').Value = '  -4.62%  '
$ws.Range('E47').Value = '  -0.40%  '
$ws.Range('D48').Value = '15.04'
$ws.Range('E48').Value = '  -6.74%  '
$ws.Range('D49').Value = '7.14'
$ws.Range('E49').Value = '  -0.19%  '
$ws.Range('E50').Value = '  +0.57%  '
$ws.Range('D51').Value = '2.252.12'
$ws.Range('E51').Value = '  -0.56%  '
